$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-16 share the same structure:
#   D = Fecha (date serial), L = Calidad, M = Volumen,
#   N = Precio minimo, O = Precio maximo, P = Precio promedio ponderado,
#   S = Precio $/Kg
# The data was re-shuffled across the 15 data rows (2..16); only these
# columns change per the diff, everything else stays the same.

$rows = @{
  2  = @{ D = 44559; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
  3  = @{ D = 44559; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 }
  4  = @{ D = 44532; L = "Primera"; M = 100; N = 10000; O = 10000; P = 10000; S = 5000 }
  5  = @{ D = 44532; L = "Segunda"; M = 100; N = 8000;  O = 8000;  P = 8000;  S = 4000 }
  6  = @{ D = 44602; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
  7  = @{ D = 44602; L = "Segunda"; M = 100; N = 5000;  O = 5000;  P = 5000;  S = 2500 }
  8  = @{ D = 44216; L = "Primera"; M = 200; N = 3500;  O = 4000;  P = 3750;  S = 1875 }
  9  = @{ D = 44216; L = "Segunda"; M = 100; N = 3000;  O = 3000;  P = 3000;  S = 1500 }
  10 = @{ D = 44574; L = "Primera"; M = 200; N = 7000;  O = 8000;  P = 7500;  S = 3750 }
  11 = @{ D = 44574; L = "Segunda"; M = 100; N = 6000;  O = 6000;  P = 6000;  S = 3000 }
  12 = @{ D = 44195; L = "Primera"; M = 200; N = 3000;  O = 3500;  P = 3250;  S = 1625 }
  13 = @{ D = 44195; L = "Segunda"; M = 100; N = 2500;  O = 2500;  P = 2500;  S = 1250 }
  14 = @{ D = 44609; L = "Primera"; M = 100; N = 6500;  O = 7000;  P = 6750;  S = 3375 }
  15 = @{ D = 44609; L = "Segunda"; M = 50;  N = 6000;  O = 6000;  P = 6000;  S = 3000 }
  16 = @{ D = 44617; L = "Primera"; M = 200; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
}

foreach ($r in $rows.Keys) {
  $v = $rows[$r]
  $ws.Cells.Item($r, 4).Value2  = $v.D   # D: Fecha
  $ws.Cells.Item($r, 12).Value  = $v.L   # L: Calidad
  $ws.Cells.Item($r, 13).Value2 = $v.M   # M: Volumen
  $ws.Cells.Item($r, 14).Value2 = $v.N   # N: Precio minimo
  $ws.Cells.Item($r, 15).Value2 = $v.O   # O: Precio maximo
  $ws.Cells.Item($r, 16).Value2 = $v.P   # P: Precio promedio ponderado
  $ws.Cells.Item($r, 19).Value2 = $v.S   # S: Precio $/Kg
}
